$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 6817.676493320005
$ws.Range("E2").Value = 202876.6317618639
$ws.Range("I2").Value = 113226.96948346
$ws.Range("L2").Value = 407171.573684766
$ws.Range("M2").Value = 74134.115776105
$ws.Range("N2").Value = 49582.36890570669
$ws.Range("O2").Value = 48727.30956737027

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 33170.24475158019
$ws.Range("E2").Value = 189865.7594528699
$ws.Range("I2").Value = 196298.3198217027
$ws.Range("L2").Value = 154748.0785685412
$ws.Range("M2").Value = 79755.84585322202
$ws.Range("N2").Value = 23752.27722812035
$ws.Range("O2").Value = 35339.83062727444

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 28619.61401238371
$ws.Range("B2").Value = 23143.29485244409
$ws.Range("E2").Value = 111916.8406725409
$ws.Range("I2").Value = 150385.2728707001
$ws.Range("M2").Value = 34803.41203795493
$ws.Range("N2").Value = 44938.11408779013
$ws.Range("O2").Value = 26938.31306104351

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 1014.766490779938

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 34409.11717595647
$ws.Range("N2").Value = 5182.698656944208
$ws.Range("O2").Value = 22972.54525065906
